$d = $word.ActiveDocument

# --- Merge split runs in title/author/abstract paragraphs into single runs ---
$d.Content.Find.Execute("Questions: The scalar product", $false, $false, $false, $false, $false, $true, 1, $false, "Questions: The scalar product", 2) | Out-Null
$d.Content.Find.Execute("Ritwik Anand", $false, $false, $false, $false, $false, $true, 1, $false, "Ritwik Anand", 2) | Out-Null
$d.Content.Find.Execute("A selection of questions for the study guide on the scalar product", $false, $false, $false, $false, $false, $true, 1, $false, "A selection of questions for the study guide on the scalar product", 2) | Out-Null

# --- Normalize m:dPr child order (begChr, sepChr, endChr) to canonical schema order in the 20 affected paragraphs ---
$d.Paragraphs(8).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">1.1.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>6</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>3</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>4</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>4</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>2</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(9).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">1.2.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>10</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>7</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>4</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>3</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>13</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(10).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">1.3.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>44</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>12</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>3</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>61</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>25</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>93</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(11).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">1.4.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>54</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>38</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>0</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>32</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>55</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>13</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(19).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">2.1.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>2</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>2</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>11</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(20).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">2.2.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(21).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">2.3.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>8</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>7</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(22).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">2.4.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>1.2</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1.4</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3.1</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5.4</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>9.7</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>7.5</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(23).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">2.5.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>45</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>65</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>54</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>19</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>58</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>71</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(24).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">2.6.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>0</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>0</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>0</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>0</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(25).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">2.7.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>3</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>4</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>6</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(26).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">2.8.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>17</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>3</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>8</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>12</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>19</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>16</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(29).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">3.1.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>2</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>4</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>7</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>λ</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(30).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">3.2.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>0</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>λ</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>2</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>3</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(31).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">3.3.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>9</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>11</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>λ</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>λ</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>3</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(32).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">3.4.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>λ</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>6</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>λ</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>λ</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>8</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(33).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">3.5.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r><m:sSup><m:e><m:r><m:t>λ</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:e></m:mr><m:mr><m:e><m:r><m:t>4</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>14</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>3</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>λ</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(34).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">3.6.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>9</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>λ</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>λ</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>λ</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(35).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">3.7.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>7</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>4</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>λ</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>λ</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>6</m:t></m:r><m:r><m:t>λ</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
$d.Paragraphs(36).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:t xml:space="preserve">3.8.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>25</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:sSup><m:e><m:r><m:t>λ</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>λ</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>11</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>7</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
